$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, shifting existing rows 214-306 down to 215-307
$ws.Rows("214").Insert()

# Populate the newly inserted row 214 with the new data record
$ws.Range("A214").Value = 10
$ws.Range("B214").Value = "Vega Modelo de Temuco"
$ws.Range("C214").Value = "La Araucanía"
$ws.Range("D214").Value = 44704
$ws.Range("E214").Value = 9
$ws.Range("F214").Value = 100114013
$ws.Range("G214").Value = "Zanahoria"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 100
$ws.Range("K214").Value = 6000
$ws.Range("L214").Value = 6000
$ws.Range("M214").Value = 6000
$ws.Range("N214").Value = "`$/saco 25 kilos"
$ws.Range("O214").Value = "Región de La Araucanía"
$ws.Range("P214").Value = 240
$ws.Range("Q214").Value = 25
$ws.Range("R214").Value = "Hortaliza"
